$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2-16 from 45207 to 45208 (+1 day)
$ws.Range("C2:C16").Value = 45208
